$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 1. Status text update: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (appears on Overview!E2/F2/E3/F3 and on the zh-cn/de-de sheets' Status col C2/C3)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Populate the "Latest Target File" (I) / "Latest Handback File" (J) /
#    "Latest Handback DateTime" (K) columns on the zh-cn and de-de sheets,
#    which previously were blank / sentinel-dated, now that handback has run.
# ---------------------------------------------------------------------------

# -- zh-cn sheet --
$zhcnTarget1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e33cb351585ead2ebb1cd6a3396bac73f18e5e8/e2e/562117bb-5af2-40a6-8713-b398be52c1ed.md"
$zhcnTarget2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e33cb351585ead2ebb1cd6a3396bac73f18e5e8/e2e/5eb3d5b0-c04d-4b6b-83f1-cd1ac62f5fc6.md"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $zhcnTarget1, "", "", "562117bb-5af2-40a6-8713-b398be52c1ed.md")
$wsZhCn.Range("I2").Style = "HyperLink"
$wsZhCn.Range("J2").Value = "562117bb-5af2-40a6-8713-b398be52c1ed.2b5052ed77dfd99f56a43e5d7fca941fdfe0267f.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-06 05:18:30"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $zhcnTarget2, "", "", "5eb3d5b0-c04d-4b6b-83f1-cd1ac62f5fc6.md")
$wsZhCn.Range("I3").Style = "HyperLink"
$wsZhCn.Range("J3").Value = "5eb3d5b0-c04d-4b6b-83f1-cd1ac62f5fc6.2b21bca05d1e17f1b085101678fe17e3ecccf175.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-06 05:18:30"

# -- de-de sheet --
$dedeTarget1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e33cb351585ead2ebb1cd6a3396bac73f18e5e8/e2e/562117bb-5af2-40a6-8713-b398be52c1ed.md"
$dedeTarget2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e33cb351585ead2ebb1cd6a3396bac73f18e5e8/e2e/5eb3d5b0-c04d-4b6b-83f1-cd1ac62f5fc6.md"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $dedeTarget1, "", "", "562117bb-5af2-40a6-8713-b398be52c1ed.md")
$wsDeDe.Range("I2").Style = "HyperLink"
$wsDeDe.Range("J2").Value = "562117bb-5af2-40a6-8713-b398be52c1ed.2b5052ed77dfd99f56a43e5d7fca941fdfe0267f.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-06 05:18:38"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $dedeTarget2, "", "", "5eb3d5b0-c04d-4b6b-83f1-cd1ac62f5fc6.md")
$wsDeDe.Range("I3").Style = "HyperLink"
$wsDeDe.Range("J3").Value = "5eb3d5b0-c04d-4b6b-83f1-cd1ac62f5fc6.2b21bca05d1e17f1b085101678fe17e3ecccf175.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-06 05:18:38"

# ---------------------------------------------------------------------------
# 3. Widen columns that now hold the longer hyperlink / file-name text so the
#    extra "Latest Target File" / "Latest Handback File" columns match the
#    other wide (40-char) columns, and the "Status" columns grow to fit the
#    new, longer status message.
# ---------------------------------------------------------------------------

# Overview: Status columns E (zh-cn) and F (de-de)
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

foreach ($ws in @($wsZhCn, $wsDeDe)) {
    # Status column C
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    # Latest Target File (I) / Latest Handback File (J) columns
    $ws.Columns.Item(9).ColumnWidth  = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}
